# Update scripts with new TPM values.
#
# The sheet lists ligand/receptor edge-weight statistics for every
# (Sending cluster, Target cluster) pair.  Columns G (Ligand average
# expression value) and M (Receptor average expression value) are the
# two "raw" TPM-derived quantities; every other changed column
# (H,I,J,N,O,P,Q,R,S,T) is a deterministic function of those raw values:
#
#   H = G * 3                              (ligand total expression)
#   I = J = G / sum(G over sending clusters)
#   N = M * 3                              (receptor total expression)
#   O = P = M / sum(M over target clusters)
#   Q = G * M                              (edge average expression weight)
#   R = H * N                              (edge total expression weight)
#   S = Q / sum(Q over all rows)
#   T = R / sum(R over all rows)
#
# The underlying TPM recompute changed the ligand/receptor averages for
# the ECs, MuSCs and Resolving-Mac clusters (FAPs is unchanged), so we
# only need to poke the new G/M numbers in and let the cascade formulas
# above reproduce every other touched cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Ligand average expression values (column G), keyed by sending
# cluster name (column A).
$GNew = @{
    "ECs"            = 35.73885133333334
    "FAPs"           = 1689.289306666667
    "MuSCs"          = 93.641553
    "Resolving-Mac"  = 14.34625366666667
}

# New Receptor average expression values (column M), keyed by target
# cluster name (column D).
$MNew = @{
    "ECs"            = 0.5550926666666666
    "FAPs"           = 3.387303666666666
    "MuSCs"          = 0.5311786666666667
    "Resolving-Mac"  = 0.117885
}

$sumG = ($GNew.Values | Measure-Object -Sum).Sum
$sumM = ($MNew.Values | Measure-Object -Sum).Sum

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

$rows = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value2   # column A
    $target  = $ws.Cells.Item($r, 4).Value2   # column D

    $g = $GNew[$sending]
    $m = $MNew[$target]

    $h = $g * 3
    $i = $g / $sumG
    $n = $m * 3
    $o = $m / $sumM
    $q = $g * $m
    $rr = $h * $n

    $ws.Cells.Item($r, 7).Value  = $g    # G
    $ws.Cells.Item($r, 8).Value  = $h    # H
    $ws.Cells.Item($r, 9).Value  = $i    # I
    $ws.Cells.Item($r, 10).Value = $i    # J
    $ws.Cells.Item($r, 13).Value = $m    # M
    $ws.Cells.Item($r, 14).Value = $n    # N
    $ws.Cells.Item($r, 15).Value = $o    # O
    $ws.Cells.Item($r, 16).Value = $o    # P
    $ws.Cells.Item($r, 17).Value = $q    # Q
    $ws.Cells.Item($r, 18).Value = $rr   # R

    $rows += , @($r, $q, $rr)
}

$sumQ = ($rows | ForEach-Object { $_[1] } | Measure-Object -Sum).Sum
$sumR = ($rows | ForEach-Object { $_[2] } | Measure-Object -Sum).Sum

foreach ($row in $rows) {
    $r  = $row[0]
    $q  = $row[1]
    $rr = $row[2]
    $s = $q / $sumQ
    $t = $rr / $sumR
    $ws.Cells.Item($r, 19).Value = $s   # S
    $ws.Cells.Item($r, 20).Value = $t   # T
}
